$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4: expand the "skip" note with more detail
$ws.Range("D4").Value = "skip, thin and long, east west"

# Fill in the "old district" (column C) mapping for the remaining rows
$ws.Range("C8").Value  = "N5"
$ws.Range("C9").Value  = "N7"
$ws.Range("C10").Value = "N6"
$ws.Range("C11").Value = "E7"
$ws.Range("C12").Value = "E8"
$ws.Range("C13").Value = "E5"
$ws.Range("C14").Value = "E9"
$ws.Range("C15").Value = "E2"
$ws.Range("C16").Value = "nothing"
$ws.Range("C17").Value = "S3"
$ws.Range("C18").Value = "W6"
$ws.Range("C19").Value = "W8"
$ws.Range("C20").Value = "E6"
$ws.Range("C21").Value = "W7"
$ws.Range("C25").Value = "W4"
$ws.Range("C26").Value = "W4"
$ws.Range("C27").Value = "nothing"
$ws.Range("C28").Value = "N8"
$ws.Range("C29").Value = "W8"
$ws.Range("C30").Value = "E1 & E3"
$ws.Range("C31").Value = "S2"
$ws.Range("C32").Value = "S5"
$ws.Range("C33").Value = "S8"
$ws.Range("C35").Value = "S6"
$ws.Range("C36").Value = "S6, S4"
$ws.Range("C37").Value = "S2, S4"
$ws.Range("C38").Value = "S9"

# Mark the newly-mapped rows' progress status (column D)
$ws.Range("D7").Value  = "done"
$ws.Range("D8").Value  = "done"
$ws.Range("D9").Value  = "done"
$ws.Range("D10").Value = "done"
$ws.Range("D11").Value = "done"
$ws.Range("D12").Value = "done"
$ws.Range("D13").Value = "done"
$ws.Range("D14").Value = "done"
$ws.Range("D15").Value = "done"
$ws.Range("D16").Value = "wip"

# D2 is no longer bold
$ws.Range("D2").Font.Bold = $false

# Move the active selection to D15
$ws.Range("D15").Select()
